$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650996069779039"
$ws1.Range("B2").Value = "go_stims-1650996069738999.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960697630353.csv"
$ws1.Range("B4").Value = "go_stims-16509960697630353.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996069779039.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509960716153831"
$ws2.Range("B2").Value = "OB-16509960701113844.csv"
$ws2.Range("B3").Value = "ZB-match_1-16509960698110092.csv"
$ws2.Range("B4").Value = "OB-16509960707833881.csv"
$ws2.Range("B5").Value = "ZB-match_5-16509960699150007.csv"
$ws2.Range("B6").Value = "TB-16509960709993858.csv"
$ws2.Range("B7").Value = "TB-16509960712794287.csv"
$ws2.Range("B8").Value = "OB-16509960702393897.csv"
$ws2.Range("B9").Value = "TB-1650996071591418.csv"
$ws2.Range("B10").Value = "ZB-match_0-16509960699991362.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509960716153831"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509960716634183"
$ws4.Range("B2").Value = "MM_stims-16509960716314187.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960716153831.csv"
$ws4.Range("B4").Value = "MM_stims-16509960716473927.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960716314187.csv"
$ws4.Range("B6").Value = "MM_stims-16509960716634183.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960716473927.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509960717273836"
$ws5.Range("B2").Value = "SAT_stims-16509960716794236.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509960717114193.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509960716954198.csv"
$ws5.Range("B5").Value = "SAT_stims-16509960716634183.csv"
